# Update course excel file: EDISON SCHOOL OF TECH SCIENCES -> per-category department
# labels, simplify the NSW/QLD location text (move the "currently not accepting
# enrolments" caveat into its own locationDetail column) and drop the expired
# promotion-validity note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("department") used to just repeat the school name for every row.
# Replace it with the actual course category for each group of rows.
$ws.Range("C2:C4").Value = "Air-Conditioning"
$ws.Range("C5:C8").Value = "Automotive"
$ws.Range("C9:C11").Value = "Packages"

# Column M ("location") drop the "(Currently not accepting enrolments)" suffix,
# keep the plain state code ...
$ws.Range("M2").Value = "NSW/QLD"
$ws.Range("M3").Value = "NSW/QLD"
$ws.Range("M4").Value = "NSW/QLD "
$ws.Range("M9").Value = "NSW/QLD"

# ... and record that caveat in the locationDetail column (N) instead.
$ws.Range("N2").Value = "Currently not accepting enrolments"
$ws.Range("N3").Value = "Currently not accepting enrolments"
$ws.Range("N4").Value = "Currently not accepting enrolments"
$ws.Range("N9").Value = "Currently not accepting enrolments"

# Column R ("promotionValidity") - the promotion has expired, clear it for all rows.
$ws.Range("R2:R11").ClearContents()
